# edit.ps1 -- PowerPoint COM-interop script reproducing the target commit.
#
# The commit does two things to the deck:
#
#   1. The table on slide 6 ("Sources of finance") is switched to a
#      different built-in table style
#      ({3E318380-...} -> {FD231CB0-62CB-4E8C-829F-7D2716B621AB}).
#
#   2. The presentation's colour theme is switched from the "Integral"
#      palette to the "Office Theme" palette (the font scheme and the
#      format/effect scheme are byte-identical between the two themes,
#      so only the 12 theme colours actually change).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 --------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FD231CB0-62CB-4E8C-829F-7D2716B621AB}")
    }
}

# --- 2) Swap the theme colours: Integral -> Office Theme -------------------
# Order matches the standard clrScheme child order / ThemeColorScheme index:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Values are Office Theme's srgbClr hex values, re-packed as COM RGB
# (0x00BBGGRR) integers for the .RGB property.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$officeColorsBGR = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColorsBGR[$i - 1]
}
